$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.878.08'
$ws.Range('E2').Value = '  +1.06%  '
$ws.Range('D3').Value = '2.466.13'
$ws.Range('E3').Value = '  +1.52%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '570.80'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.25%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.62'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.44%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('E8').Value = '  +1.04%  '
$ws.Range('D9').Value = '2.466.06'
$ws.Range('E9').Value = '  +1.50%  '
$ws.Range('E10').Value = '  +1.35%  '
$ws.Range('E11').Value = '  +0.82%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.24'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.19%  '
$ws.Range('E13').Value = '  +1.67%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.09'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.88%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000178'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.61%  '
$ws.Range('E16').Value = '  +1.62%  '
$ws.Range('D17').Value = '62.833.01'
$ws.Range('E17').Value = '  +1.09%  '
$ws.Range('D18').Value = '2.466.53'
$ws.Range('E18').Value = '  +1.42%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.31'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +6.67%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '323.34'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.31%  '
$ws.Range('E22').Value = '  +0.36%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.93'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +11.36%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.998'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.18%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '65.78'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.55%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '618.52'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +11.31%  '
$ws.Range('E27').Value = '  +8.48%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.49'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.81%  '
$ws.Range('D29').Value = '2.589.67'
$ws.Range('E29').Value = '  +1.77%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.999'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.08%  '
$ws.Range('E31').Value = '  +4.68%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.18'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.21%  '
$ws.Range('B33').Value = 'PancakeSwap'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.90'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.74%  '
$ws.Range('B34').Value = 'Kaspa'
$ws.Range('C34').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.141'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.31%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.06'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.61%  '
$ws.Range('E36').Value = '  -2.76%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.11%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.382'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.16%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.39'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.94%  '
$ws.Range('E40').Value = '  +0.05%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '145.80'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.35%  '
$ws.Range('E42').Value = '  -0.95%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.57'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +12.88%  '
$ws.Range('E44').Value = '  +0.02%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '147.76'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.44%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.72'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.57%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '20.72'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.53%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0538'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.86%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.604'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.20%  '
$ws.Range('E50').Value = '  +1.23%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0917'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.97%  '
